# Rename the "investment" category items to reflect the move to the
# inline bot: prefix each of the investment rows (A38:A43) with the new
# Hebrew label text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = "השקעה - שכירות"
$ws.Range("A39").Value = "השקעה - ביטוח דירה"
$ws.Range("A40").Value = "השקעה - תחזוקה"
$ws.Range("A41").Value = "השקעה - ניהול"
$ws.Range("A42").Value = "השקעה - ריהוט"
$ws.Range("A43").Value = "השקה - משכנתא"

# Restore the saved selection / scroll position recorded for this sheet.
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E37").Select()
